$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "35.002.38"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.32%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.848.07"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +2.26%  "
$ws.Range("E4").Value = "  +0.30%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "232.82"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.618"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("E7").Value = "  +0.36%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "40.87"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +5.20%  "
$ws.Range("E9").Value = "  +3.53%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0692"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +2.03%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0987"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.60%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "2.119.07"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.46%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "11.42"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +5.10%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "1.849.84"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.90%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.675"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.63%  "
$ws.Range("E16").Value = "  +2.56%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "35.061.44"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.67%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "70.10"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.52%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.0₃0791"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.63%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "240.65"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.05%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.26"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +4.08%  "
$ws.Range("E22").Value = "  +3.02%  "
$ws.Range("E23").Value = "  +0.31%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.26"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.35%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "172.68"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.32%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "7.86"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.51%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "17.52"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("E28").Value = "  +3.94%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.60"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +3.51%  "
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("E32").Value = "  -0.35%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.97"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.19%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.60"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +22.53%  "
$ws.Range("E35").Value = "  +11.82%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.754"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +10.78%  "
$ws.Range("E37").Value = "  +7.33%  "
$ws.Range("E38").Value = "  +12.92%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "90.30"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.34%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.351.09"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.42%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0197"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +3.08%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "14.73"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.28%  "
$ws.Range("E44").Value = "  -1.64%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.76"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("E46").Value = "  +4.00%  "
$ws.Range("E47").Value = "  +3.34%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.038.02"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.53%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "3.41"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +18.84%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0673"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.08%  "

Write-Host "Applied all changes"
